$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad" / last-changed date) holds serial date 45207 (2023-10-08)
# for every data row (rows 2 through 157). The update bumps it by one day to
# 45208 (2023-10-09) for all of those rows.
for ($r = 2; $r -le 157; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value = 45208
    }
}

# Row 2's hyperlink formulas (columns S, T, V, W, X, Y) pointed at the
# "Logging_SKELLEFTEA" folder; they now point at "Logging_2482" instead.
$hyperlinkCols = @(19, 20, 22, 23, 24, 25)  # S, T, V, W, X, Y
foreach ($col in $hyperlinkCols) {
    $cell = $ws.Cells.Item(2, $col)
    $formula = $cell.Formula
    if ($formula -and $formula -like "*Logging_SKELLEFTEA*") {
        $cell.Formula = $formula -replace "Logging_SKELLEFTEA", "Logging_2482"
    }
}
